# Update countries & provincias Spain
# Applies the data refresh for the "Pais" worksheet:
#  - Updates the "Datos actualizados" timestamp in A1
#  - Updates several countries' statistics (columns B-H)
#  - Malta overtakes Nueva Zelanda in total cases -> rows 143/144 swap countries
#  - Montserrat overtakes Islas Malvinas in total cases -> rows 214/215 swap countries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp update
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 13:51"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5916089
$ws.Range("C4").Value = 459
$ws.Range("E4").Value = 2516458

# Kuwait (row 41)
$ws.Range("B41").Value = 81573
$ws.Range("C41").Value = 613
$ws.Range("D41").Value = 73402
$ws.Range("E41").Value = 7652
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 519

# Suiza (row 61)
$ws.Range("B61").Value = 40262
$ws.Range("C61").Value = 202
$ws.Range("E61").Value = 3861

# Bosnia y Herzegovina (row 77)
$ws.Range("B77").Value = 18326
$ws.Range("C77").Value = 297
$ws.Range("D77").Value = 12081
$ws.Range("E77").Value = 5685
$ws.Range("G77").Value = 13
$ws.Range("H77").Value = 560

# Madagascar (row 82)
$ws.Range("B82").Value = 14475
$ws.Range("C82").Value = 73
$ws.Range("D82").Value = 13492
$ws.Range("E82").Value = 805

# Malta / Nueva Zelanda swap (rows 143-144)
# Row 143 becomes Malta with updated stats
$ws.Range("A143").Value = "Malta"
$ws.Range("B143").Value = 1705
$ws.Range("C143").Value = 38
$ws.Range("D143").Value = 1029
$ws.Range("E143").Value = 666
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 10

# Row 144 becomes Nueva Zelanda with its former (unchanged) stats
$ws.Range("A144").Value = "Nueva Zelanda"
$ws.Range("B144").Value = 1690
$ws.Range("C144").Value = 7
$ws.Range("D144").Value = 1539
$ws.Range("E144").Value = 129
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 22

# Vietnam (row 161)
$ws.Range("B161").Value = 1028
$ws.Range("C161").Value = 6
$ws.Range("E161").Value = 411

# Liechtenstein (row 194)
$ws.Range("B194").Value = 102
$ws.Range("C194").Value = 2
$ws.Range("D194").Value = 94
$ws.Range("E194").Value = 7

# Montserrat / Islas Malvinas swap (rows 214-215)
# Row 214 becomes Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

# Row 215 becomes Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
